$wb = $excel.ActiveWorkbook

# Rename the "Include from SmokingStatus" sheet to "Include from Smoking Status"
$includeSheet = $wb.Worksheets.Item("Include from SmokingStatus")
$includeSheet.Name = "Include from Smoking Status"

# Metadata sheet updates
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.0.0"
$meta.Range("B5").Value = "Smoking Status"

# Experimental value needs to land as literal text "false" (not a Boolean),
# so enter it as a formula producing the text and paste back as a value.
$expCell = $meta.Range("B7")
$expCell.Formula = "=""false"""
$expCell.Copy()
$expCell.PasteSpecial(-4163)

$meta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Include sheet: update the System URI value
$includeSheet.Range("B9").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/SmokingStatusCS"
